# Apply the edit described by the diff:
#  - the title (column A) and uri (column E) text for row 2 and row 3 swap
#    with each other, and likewise for row 4 and row 5 (net effect of the
#    shared-string table being re-serialized after a json was added to the
#    backing dataset for the time-bucket analysis).
#  - the hyperlink "ref" and its r:id stay pinned to the same cell (E2 always
#    uses rId1, E3 always uses rId2, etc.) but the location="1960"
#    sub-address moves from the E2 hyperlink to the E3 hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- swap the title text (column A) between row 2 and row 3 ---
$a2 = $ws.Range("A2").Value()
$a3 = $ws.Range("A3").Value()
$ws.Range("A2").Value = $a3
$ws.Range("A3").Value = $a2

# --- swap the title text (column A) between row 4 and row 5 ---
$a4 = $ws.Range("A4").Value()
$a5 = $ws.Range("A5").Value()
$ws.Range("A4").Value = $a5
$ws.Range("A5").Value = $a4

# --- swap the uri text (column E) between row 2 and row 3 ---
$e2 = $ws.Range("E2").Value()
$e3 = $ws.Range("E3").Value()
$ws.Range("E2").Value = $e3
$ws.Range("E3").Value = $e2

# --- swap the uri text (column E) between row 4 and row 5 ---
$e4 = $ws.Range("E4").Value()
$e5 = $ws.Range("E5").Value()
$ws.Range("E4").Value = $e5
$ws.Range("E5").Value = $e4

# --- move the location="1960" sub-address from the E2 hyperlink to the E3
#     hyperlink, leaving each hyperlink's ref/r:id (and therefore its
#     Address) untouched ---
$hlArr = @($ws.Hyperlinks)
$hlArr[0].SubAddress = ""
$hlArr[1].SubAddress = "1960"
